# Update column F (dSF) values to match repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -2
    13 = -6
    15 = 8
    17 = -5
    22 = -6
    23 = -2
    25 = 0
    28 = 1
    31 = -4
    34 = -6
    36 = -3
    37 = 6
    38 = -3
    39 = -12
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
